$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking price strings so they are not
# auto-converted to numbers by Excel (matches original inline-string text cells).
$textCells = @('D5','D6','D7','D9','D10','D12','D14','D15','D16','D20','D22','D25','D26','D27','D33','D34','D35','D36','D37','D38','D39','D40','D42','D46','D47','D48','D49')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '39.480.50'
$ws.Range('E2').Value = '  +1.89%  '
$ws.Range('D3').Value = '2.164.56'
$ws.Range('E3').Value = '  +2.93%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '228.13'
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('D6').Value = '0.624'
$ws.Range('E6').Value = '  +1.16%  '
$ws.Range('D7').Value = '63.95'
$ws.Range('E7').Value = '  +2.91%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').Value = '0.396'
$ws.Range('E9').Value = '  +1.90%  '
$ws.Range('D10').Value = '0.0856'
$ws.Range('E10').Value = '  +1.38%  '
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('D12').Value = '16.10'
$ws.Range('E12').Value = '  +1.89%  '
$ws.Range('D13').Value = '2.484.29'
$ws.Range('E13').Value = '  +2.95%  '
$ws.Range('D14').Value = '22.16'
$ws.Range('E14').Value = '  +0.23%  '
$ws.Range('D15').Value = '0.815'
$ws.Range('E15').Value = '  +0.97%  '
$ws.Range('D16').Value = '5.54'
$ws.Range('E16').Value = '  +0.44%  '
$ws.Range('D17').Value = '2.166.16'
$ws.Range('E17').Value = '  +2.79%  '
$ws.Range('D18').Value = '39.471.07'
$ws.Range('E18').Value = '  +1.82%  '
$ws.Range('E19').Value = '  +0.98%  '
$ws.Range('D20').Value = '71.92'
$ws.Range('E20').Value = '  -0.01%  '
$ws.Range('D21').Value = '0.0₃0852'
$ws.Range('E21').Value = '  +1.28%  '
$ws.Range('D22').Value = '229.74'
$ws.Range('E22').Value = '  +0.70%  '
$ws.Range('E24').Value = '  -0.88%  '
$ws.Range('D25').Value = '2.37'
$ws.Range('E25').Value = '  +1.75%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '172.30'
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = '9.52'
$ws.Range('E27').Value = '  -0.69%  '
$ws.Range('E28').Value = '  +2.33%  '
$ws.Range('E29').Value = '  +2.71%  '
$ws.Range('E30').Value = '  -0.38%  '
$ws.Range('E31').Value = '  +4.26%  '
$ws.Range('E32').Value = '  +1.18%  '
$ws.Range('D33').Value = '4.63'
$ws.Range('E33').Value = '  +1.83%  '
$ws.Range('D34').Value = '7.09'
$ws.Range('E34').Value = '  +1.87%  '
$ws.Range('D35').Value = '4.73'
$ws.Range('E35').Value = '  -0.53%  '
$ws.Range('D36').Value = '0.0621'
$ws.Range('E36').Value = '  +0.40%  '
$ws.Range('D37').Value = '2.45'
$ws.Range('E37').Value = '  +0.80%  '
$ws.Range('D38').Value = '3.59'
$ws.Range('E38').Value = '  -0.80%  '
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  +0.06%  '
$ws.Range('D40').Value = '103.23'
$ws.Range('E40').Value = '  +0.22%  '
$ws.Range('E41').Value = '  +0.52%  '
$ws.Range('D42').Value = '17.86'
$ws.Range('E42').Value = '  -1.57%  '
$ws.Range('D43').Value = '1.523.78'
$ws.Range('E44').Value = '  +3.38%  '
$ws.Range('E45').Value = '  +5.69%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').Value = '0.0928'
$ws.Range('E46').Value = '  +1.79%  '
$ws.Range('B47').Value = 'HuobiToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D47').Value = '2.82'
$ws.Range('E47').Value = '  +0.74%  '
$ws.Range('D48').Value = '4.26'
$ws.Range('E48').Value = '  +3.37%  '
$ws.Range('D49').Value = '7.75'
$ws.Range('E49').Value = '  -1.09%  '
$ws.Range('D50').Value = '2.367.56'
$ws.Range('E50').Value = '  +2.95%  '
$ws.Range('E51').Value = '  -0.55%  '
